# From v1.2 to v1.2.1
# The TC4/TC5 "Steps"/"Expected Results" content was swapped:
#   - TC4 now shows the "detalhar a solicitação de diária" step/result
#   - TC5 now shows the "excluir comprovante" step/result
# And the TC7/TC8 "Expected Results" content was swapped:
#   - TC7 now shows the MSG212 (ainda não pode ser realizada) message
#   - TC8 now shows the "não está em nenhum desses dois estados" message

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- TC4 / TC5 swap (Steps column B, Expected Results column D) ---
$tc4Step = $ws.Range("B41").Value2
$tc4Result = $ws.Range("D41").Value2
$tc5Step = $ws.Range("B50").Value2
$tc5Result = $ws.Range("D50").Value2

$ws.Range("B41").Value = $tc5Step
$ws.Range("D41").Value = $tc5Result
$ws.Range("B50").Value = $tc4Step
$ws.Range("D50").Value = $tc4Result

# --- TC7 / TC8 swap (Expected Results column D) ---
$tc7Result = $ws.Range("D67").Value2
$tc8Result = $ws.Range("D75").Value2

$ws.Range("D67").Value = $tc8Result
$ws.Range("D75").Value = $tc7Result
